$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Company" column (column C) entirely - shifts State/Enquiry For/
# Enquiry From left by one column.
$ws.Columns("C").Delete()

# Update header: "Enquiry For" -> "Template Name"
$ws.Range("D1").Value = "Template Name"

# Update row 2 values
$ws.Range("C2").Value = "Rajasthan"
$ws.Range("D2").Value = "LAVANYA_HINDI"
$ws.Range("E2").Value = "Indiamart"

# Update row 3 values
$ws.Range("D3").Value = "LAVANYA_ENGLISH"
$ws.Range("E3").Value = "Exporters India"
